# Add a new "2022-Q3" sheet (right after the "总计" summary sheet) and
# insert a matching summary row on "总计" for the new quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (totals) sheet: insert a new row 2 for 2022-Q3.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Rows.Item(2).Insert()

# Copy the formatting (bold / centered / bordered) used by the other
# index cells in column A onto the freshly inserted A2 cell.
$total.Cells.Item(3, 1).Copy()
$total.Cells.Item(2, 1).PasteSpecial(-4122)   # xlPasteFormats
$total.Range("B2:D2").ClearFormats()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 10
$total.Cells.Item(2, 4).Value = 1.3

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet right after "总计" and fill it
#    with the quarter's fund-holding detail table.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")

$headerRange = $q3.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
for ($col = 2; $col -le 8; $col++) {
    $q3.Cells.Item(1, $col).Value = $headers[$col - 2]
}

$indexRange = $q3.Range("A2:A11")
$indexRange.Font.Bold = $true
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160
$indexRange.Borders.LineStyle = 1

# Numeric-looking columns must be forced to text so values such as
# "050001" / "21.72" keep leading zeros / exact formatting instead of
# being auto-converted to numbers by Excel's smart entry. (Row 11's G
# cell is a genuine number in the source data, so G is only forced to
# text through row 10.)
$q3.Range("B2:B11").NumberFormat = "@"
$q3.Range("D2:F11").NumberFormat = "@"
$q3.Range("G2:G10").NumberFormat = "@"

$rows = @(
    @(2, "050001", "博时价值增长混合", "21.72", "72.10", "2.73", "0.5930", 8),
    @(3, "161219", "国投瑞银新兴产业混合（LOF）", "6.18", "79.94", "5.26", "0.3251", 3),
    @(4, "050201", "博时价值增长贰号混合", "10.24", "71.73", "2.38", "0.2437", 10),
    @(5, "002628", "招商安博灵活配置混合A", "0.98", "78.15", "5.36", "0.0525", 4),
    @(6, "002629", "招商安博灵活配置混合C", "0.53", "78.15", "5.36", "0.0284", 4),
    @(7, "010503", "招商稳兴混合A", "1.22", "37.98", "2.26", "0.0276", 6),
    @(8, "519097", "新华中小市值优选混合", "0.71", "67.35", "2.93", "0.0208", 8),
    @(9, "013584", "招商品质领航混合C", "0.05", "71.52", "5.88", "0.0029", 1),
    @(10, "013583", "招商品质领航混合A", "0.03", "71.52", "5.88", "0.0018", 1),
    @(11, "010504", "招商稳兴混合C", "0.00", "37.98", "2.26", $null, 6)
)

foreach ($row in $rows) {
    $r = $row[0]
    $q3.Cells.Item($r, 1).Value = $r - 2
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    if ($row[6] -ne $null) {
        $q3.Cells.Item($r, 7).Value = $row[6]
    } else {
        # G11 is a genuine number (0) in the source data, left in the
        # default/General format.
        $q3.Cells.Item($r, 7).Value = 0
    }
    $q3.Cells.Item($r, 8).Value = $row[7]
}
